$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.267.27'
$ws.Range("E2").Value = '  -0.23%  '
$ws.Range("D3").Value = '1.907.88'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '0.725'
$ws.Range("E5").Value = '  +9.83%  '
$ws.Range("D6").Value = '253.19'
$ws.Range("E6").Value = '  +3.25%  '
$ws.Range("D8").Value = '40.51'
$ws.Range("E8").Value = '  -1.91%  '
$ws.Range("E9").Value = '  +3.47%  '
$ws.Range("E10").Value = '  -1.14%  '
$ws.Range("D11").Value = '0.0761'
$ws.Range("E11").Value = '  +6.13%  '
$ws.Range("D12").Value = '0.0990'
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").Value = '2.182.92'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("E14").Value = '  +5.26%  '
$ws.Range("E15").Value = '  +1.95%  '
$ws.Range("D16").Value = '1.953.64'
$ws.Range("E16").Value = '  +2.31%  '
$ws.Range("D17").Value = '4.90'
$ws.Range("E17").Value = '  +0.65%  '
$ws.Range("D18").Value = '35.263.66'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("D19").Value = '74.26'
$ws.Range("E19").Value = '  +3.02%  '
$ws.Range("D20").Value = '0.0₃0845'
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("D21").Value = '243.03'
$ws.Range("E21").Value = '  +1.47%  '
$ws.Range("D22").Value = '12.99'
$ws.Range("E22").Value = '  +3.89%  '
$ws.Range("D23").Value = '5.08'
$ws.Range("E23").Value = '  +5.00%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("E25").Value = '  +3.45%  '
$ws.Range("D26").Value = '2.44'
$ws.Range("E26").Value = '  +3.40%  '
$ws.Range("D27").Value = '167.00'
$ws.Range("E27").Value = '  -1.78%  '
$ws.Range("D28").Value = '8.61'
$ws.Range("E28").Value = '  +1.35%  '
$ws.Range("D29").Value = '18.68'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("E30").Value = '  +4.39%  '
$ws.Range("D31").Value = '4.129.93'
$ws.Range("E31").Value = '  +19.49%  '
$ws.Range("E32").Value = '  +4.28%  '
$ws.Range("B33").Value = 'TrustWalletToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D33").Value = '1.64'
$ws.Range("E33").Value = '  +23.88%  '
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.99'
$ws.Range("E34").Value = '  +13.86%  '
$ws.Range("E35").Value = '  +2.37%  '
$ws.Range("E36").Value = '  +1.63%  '
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("D38").Value = '0.917'
$ws.Range("E38").Value = '  -1.94%  '
$ws.Range("E39").Value = '  -1.15%  '
$ws.Range("B40").Value = 'InjectiveProtocol'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D40").Value = '17.08'
$ws.Range("E40").Value = '  +4.76%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0216'
$ws.Range("E41").Value = '  +3.95%  '
$ws.Range("D42").Value = '96.29'
$ws.Range("E42").Value = '  +6.66%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '0.0647'
$ws.Range("E44").Value = '  -2.75%  '
$ws.Range("D45").Value = '1.336.88'
$ws.Range("E46").Value = '  +1.53%  '
$ws.Range("E47").Value = '  +0.58%  '
$ws.Range("D48").Value = '6.71'
$ws.Range("E48").Value = '  +2.16%  '
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("D50").Value = '45.45'
$ws.Range("E50").Value = '  -5.06%  '
$ws.Range("D51").Value = '11.99'
$ws.Range("E51").Value = '  +14.73%  '
